# eventsliders.xlsx maintenance edit
#
# - adds Artisan Command "keyboard" to enable/disable keyboard mode
#   (new row in the "Commands" sheet, just above the existing
#   "RC Command" section, listing the command name + its description)
# - (the remaining bullet points from the commit message - ambient
#   phidgets kept attached until app termination, dropped Probat
#   middleware support, recovered zh_TW translations - are changes to
#   the Artisan Python application itself and are not reflected in this
#   help-dialog workbook)

$wb = $excel.ActiveWorkbook

$wsSliders  = $wb.Worksheets.Item(1)   # "Sliders"
$wsCommands = $wb.Worksheets.Item(2)   # "Commands"

# --- Commands sheet: insert the new "keyboard" command row -----------------
# The existing table lists one command per row in columns B (command
# syntax) / C (description), with section header rows using column A.
# The new "keyboard(<bool>)" entry belongs right before the "RC Command"
# section header, which currently sits on row 95 - so insert a fresh row
# there and push everything else down by one.
$wsCommands.Rows.Item(95).Insert()

$wsCommands.Cells.Item(95, 2).Value = "keyboard(<bool>)"
$wsCommands.Cells.Item(95, 3).Value = "enables/disables keyboard mode"

# match the compact row height used by the other "command name / one-line
# description" rows in this table (e.g. row 94, 108-111)
$wsCommands.Rows.Item(95).RowHeight = 13.8

# --- view/selection bookkeeping (best effort) -------------------------------
# Keep the Sliders sheet's prior selection (B6) ...
$wsSliders.Activate() | Out-Null
$wsSliders.Range("B6").Select() | Out-Null

# ... and point the Commands sheet's selection/scroll position at the
# newly inserted row so it is visible.
$wsCommands.Activate() | Out-Null
$wsCommands.Rows.Item(95).Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 85
} catch {
}
